$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: hyperlink to external url, display text pre-set so no `display` attr is written ---
$ws1.Range("A4").Value = "hyperlink1"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://example.com/hyperlink1")

# --- Sheet2: hyperlink whose display text (from Add) differs from the final cell value ---
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://apple.com/", "", "", "https://apple.com/")
$ws2.Range("A4").Value = "hyperlink2"

# --- Sheet1: hyperlink whose display text equals the url itself ---
$ws1.Range("A5").Value = "https://google.com/"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://google.com/")

# --- Sheet1: internal/local hyperlink (location only), with explicit display text ---
$ws1.Hyperlinks.Add($ws1.Range("A6"), "", "Sheet1!A1", "", "Sheet1!A1")

# --- Sheet1: mailto hyperlink ---
$ws1.Range("A7").Value = "mailto:invalid@example.com?subject=important"
$ws1.Hyperlinks.Add($ws1.Range("A7"), "mailto:invalid@example.com?subject=important")

# --- Sheet1: string concatenation formula ---
$ws1.Range("C2").Value = "a"
$ws1.Range("C3").Value = "b"
$ws1.Range("C1").Formula = "=CONCATENATE(C2,C3)"

# --- Sheet2: numeric multiplication formula ---
$ws2.Range("C2").Value = 2
$ws2.Range("C3").Value = 3
$ws2.Range("C1").Formula = "=C2*C3"

# --- Selection state on both sheets ---
$ws2.Range("C1").Select() | Out-Null
$ws1.Range("C1").Select() | Out-Null
$ws1.Activate() | Out-Null
